$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint")

# Reword the task description for row 8 (shared string content change).
$ws.Range("C8").Value = "Project description documentation"

# Record 3 units of effort spent on day 7 ("M" column) for the
# "Expand timeline for tracks" task (row 11).
$ws.Range("M11").Value = 3

# Re-enter the daily-total and burndown formulas across their row so the
# recalculated values ripple through (this also causes Excel to store them
# as shared formulas, matching the source workbook's layout).
$ws.Range("G17:T17").Formula = "=SUM(G5:G16)"
$ws.Range("G19:T19").Formula = "=F19-SUM(G5:G16)"

# Leave the cursor on F10, as in the edited workbook.
$ws.Range("F10").Select()
